$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) The "_GoBack" bookmark currently sits between "verschieden" and
#    "e Zeitzonen ..." (an artifact of an earlier edit session). It needs to
#    move to a new spot later in the document (right after "Vienna", before
#    the closing "<<")" mark in the timezone example). Remove it from its
#    old location first; we'll re-add it in its new spot below.
# ---------------------------------------------------------------------------
$goBack = $d.Bookmarks.Item("_GoBack")
$splitPos = $goBack.Start
$goBack.Delete()

# Merge the two runs that used to be split by the bookmark ("verschieden" /
# "e Zeitzonen eingestellt, ...") back into a single run by deleting across
# the old boundary and retyping the joint - this collapses the run split
# without touching the neighbouring "Timestamp" run (which keeps its
# spell-check proofErr wrapping, unchanged).
$seam = $d.Range($splitPos - 1, $splitPos + 1)
$seamText = $seam.Text
$seam.Delete()
$seamIns = $d.Range($splitPos - 1, $splitPos - 1)
$seamIns.InsertAfter($seamText)

# ---------------------------------------------------------------------------
# 2) Replace the example timezone "Atlantic/Canary" (spell-check-flagged,
#    foreign words) with "Europe/Vienna" (not flagged). Remove the whole
#    bracketed example - including the spell-check proofErr markers around
#    "Atlantic" and "Canary" - and retype it fresh so no stray proofErr tags
#    remain.
# ---------------------------------------------------------------------------
$exampleRng = $d.Content
$found = $exampleRng.Find.Execute(" (" + [char]0x201E + "Atlantic/Canary" + [char]0x201C + ")", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$newTimezonePos = $exampleRng.Start
$exampleRng.Delete()

$ins = $d.Range($newTimezonePos, $newTimezonePos)
$ins.InsertAfter(" (" + [char]0x201E + "Europe/Vienna" + [char]0x201C + ")")

# ---------------------------------------------------------------------------
# 3) Re-add the "_GoBack" bookmark right after "Vienna" (before the closing
#    quote + parenthesis).
# ---------------------------------------------------------------------------
$closeMark = [char]0x201C + ")"
$bmRng = $d.Content
$bmRng.Find.Execute($closeMark, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$bmPos = $bmRng.Start
$target = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $target)

Write-Output "Para4: $($d.Paragraphs.Item(4).Range.Text)"
Write-Output "Para6: $($d.Paragraphs.Item(6).Range.Text)"
$bmCheck = $d.Bookmarks.Item("_GoBack")
Write-Output "Bookmark start/end: $($bmCheck.Start) / $($bmCheck.End)"
